$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old single department label "FACULTY OF BUSINESS & TECHNOLOGY"
# with per-row department categories in column C.
$ws.Range("C2:C12").Value = "Business"
$ws.Range("C13:C15").Value = "Information Technology"
$ws.Range("C16").Value = "Building and Construction"
$ws.Range("C17:C22").Value = "Packages"

# Clear the old promotion-validity text that used to sit in column R for
# every data row (the promotion expired / no longer applies).
$ws.Range("R2:R22").ClearContents()

# Restore the sheet's selection state.
$ws.Range("R2:R22").Select()
